$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Tomek" ranking row entirely (row 2)
$ws.Range("A2:D2").ClearContents()

# Update the remaining ranking entry (row 6) with the new player's results
$ws.Range("A6").Value = "Filip"
$ws.Range("B6").Value = "00:04"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Galactic Tower"

# Move the active selection to H6
$ws.Range("H6").Select()
